$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21 (HOUR 20)
$ws.Range("B21").Value = 1520.09097147339
$ws.Range("D21").Value = -73479.90902852661

# Row 22 (HOUR 21)
$ws.Range("B22").Value = 68523.995
$ws.Range("D22").Value = -6476.005000000005

# Row 23 (HOUR 22)
$ws.Range("B23").Value = 66852.495
$ws.Range("D23").Value = -8147.505000000005

# Row 24 (HOUR 23) - B changes and D24 is newly added
$ws.Range("B24").Value = 33754.2325
$ws.Range("D24").Value = -38245.7675
